$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format before writing values, so that
# numeric-looking strings (e.g. "1.001", "25.866.92") are preserved exactly
# as text instead of being auto-converted to numbers.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '25.866.92'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.630.75'
$ws.Range("E3").Value = '  -0.68%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  -0.48%  '
$ws.Range("D5").Value = '215.57'
$ws.Range("E5").Value = '  +0.20%  '
$ws.Range("D6").Value = '0.5071'
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("E7").Value = '  -0.34%  '
$ws.Range("D8").Value = '0.2574'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").Value = '0.06323'
$ws.Range("E9").Value = '  -0.95%  '
$ws.Range("D10").Value = '19.48'
$ws.Range("E10").Value = '  -0.09%  '
$ws.Range("D11").Value = '0.07752'
$ws.Range("E11").Value = '  -0.32%  '
$ws.Range("D12").Value = '4.247'
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '1.632.49'
$ws.Range("E13").Value = '  -1.11%  '
$ws.Range("D14").Value = '1.851.70'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").Value = '0.5486'
$ws.Range("E15").Value = '  +0.79%  '
$ws.Range("D16").Value = '63.66'
$ws.Range("E16").Value = '  -1.03%  '
$ws.Range("D17").Value = '0.0₅7641'
$ws.Range("E17").Value = '  -2.29%  '
$ws.Range("D18").Value = '25.874.97'
$ws.Range("E18").Value = '  -0.29%  '
$ws.Range("E19").Value = '  -0.42%  '
$ws.Range("D20").Value = '4.417'
$ws.Range("E20").Value = '  -0.59%  '
$ws.Range("D21").Value = '194.27'
$ws.Range("E21").Value = '  -1.37%  '
$ws.Range("D22").Value = '9.870'
$ws.Range("E22").Value = '  -0.68%  '
$ws.Range("D23").Value = '6.033'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  -0.53%  '
$ws.Range("D25").Value = '1.915'
$ws.Range("D26").Value = '142.02'
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").Value = '0.1242'
$ws.Range("E27").Value = '  +3.54%  '
$ws.Range("D28").Value = '6.771'
$ws.Range("E28").Value = '  -1.28%  '
$ws.Range("E29").Value = '  -0.85%  '
$ws.Range("D30").Value = '1.237'
$ws.Range("E30").Value = '  -0.03%  '
$ws.Range("D31").Value = '0.04873'
$ws.Range("E31").Value = '  -1.53%  '
$ws.Range("E32").Value = '  -0.42%  '
$ws.Range("D33").Value = '3.189'
$ws.Range("D34").Value = '1.543'
$ws.Range("E34").Value = '  +0.23%  '
$ws.Range("D35").Value = '2.372'
$ws.Range("E35").Value = '  -0.16%  '
$ws.Range("D36").Value = '0.8948'
$ws.Range("E36").Value = '  +0.08%  '
$ws.Range("D37").Value = '0.5517'
$ws.Range("E37").Value = '  +1.82%  '
$ws.Range("D38").Value = '2.538'
$ws.Range("E38").Value = '  -1.76%  '
$ws.Range("D39").Value = '1.124.31'
$ws.Range("E39").Value = '  -1.27%  '
$ws.Range("D40").Value = '0.01549'
$ws.Range("E40").Value = '  -0.24%  '
$ws.Range("E41").Value = '  -0.41%  '
$ws.Range("D42").Value = '5.569'
$ws.Range("E42").Value = '  -0.19%  '
$ws.Range("D43").Value = '0.7968'
$ws.Range("E43").Value = '  -2.24%  '
$ws.Range("E44").Value = '  -2.42%  '
$ws.Range("B45").Value = 'BabyDogeCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D45").Value = '0.0₈120'
$ws.Range("E45").Value = '  -4.37%  '
$ws.Range("B46").Value = 'RocketPoolETH'
$ws.Range("C46").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D46").Value = '1.764.96'
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("D47").Value = '0.4444'
$ws.Range("E47").Value = '  -2.35%  '
$ws.Range("D48").Value = '1.002'
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("D49").Value = '54.76'
$ws.Range("E49").Value = '  -0.05%  '
$ws.Range("D50").Value = '0.05134'
$ws.Range("E50").Value = '  +1.24%  '
$ws.Range("D51").Value = '7.540'
$ws.Range("E51").Value = '  +2.40%  '

# Restore the default (Normal) style on column D now that the text values
# have been written, so no stray number-format styling is left behind.
$dRange.Style = "Normal"
